# New crime data collected - weekly CompStat update
# (Volume/date header bump + refreshed weekly/28-day/YTD crime figures)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-TextCell($ws, $addr, $styleSourceAddr, $text) {
    # Clone the number-format/style of a cell that is already a "text style"
    # cell (s=14, General numFmt) onto the target, then force the value to be
    # stored as a genuine text value (not auto-coerced to a number) by
    # temporarily switching to the Text number format, and finally re-apply
    # the cloned General-format style so the cell matches the text-style
    # cells used elsewhere on the sheet (e.g. "0" / "***.*" placeholders).
    $ws.Range($styleSourceAddr).Copy()
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($styleSourceAddr).Copy()
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}

function Set-NumberCell($ws, $addr, $styleSourceAddr, $num) {
    # Clone the style of a cell that already carries the numeric style
    # (s=15, #,##0) onto the target, then write the number.
    $ws.Range($styleSourceAddr).Copy()
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
    $ws.Range($addr).Value = $num
}

# ---------------------------------------------------------------------------
# Report header: volume number 47 -> 48, week covered 11/20-11/26 -> 11/27-12/3
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/27/2023  Through  12/3/2023"

# ---------------------------------------------------------------------------
# Row 15 - Rape: 28-day 2022 count (G) & %chg (H) go from real numbers to the
# "no prior incidents" text placeholders ("0" / "***.*")
# ---------------------------------------------------------------------------
Set-TextCell $ws "G15" "F15" "0"
Set-TextCell $ws "H15" "F15" "***.*"

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
Set-TextCell $ws "D16" "C14" "0"
Set-TextCell $ws "E16" "E14" "***.*"
$ws.Range("F16").Value = 3
$ws.Range("H16").Value = 200
$ws.Range("I16").Value = 20
$ws.Range("K16").Value = 11.111111111111
$ws.Range("L16").Value = 100
$ws.Range("M16").Value = -9.090909090909
$ws.Range("N16").Value = -67.741935483871

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 4
$ws.Range("H17").Value = 100
$ws.Range("I17").Value = 84
$ws.Range("J17").Value = 47
$ws.Range("K17").Value = 78.723404255319
$ws.Range("L17").Value = 180
$ws.Range("M17").Value = 82.608695652173
$ws.Range("N17").Value = -17.647058823529

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -44.444444444444
$ws.Range("I18").Value = 52
$ws.Range("J18").Value = 41
$ws.Range("K18").Value = 26.829268292682
$ws.Range("L18").Value = 126.086956521739
$ws.Range("M18").Value = -47.474747474747
$ws.Range("N18").Value = -83.544303797468

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 4
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 18
$ws.Range("H19").Value = 27.777777777777
$ws.Range("I19").Value = 265
$ws.Range("J19").Value = 250
$ws.Range("K19").Value = 6
$ws.Range("L19").Value = 48.876404494382
$ws.Range("M19").Value = 93.430656934306
$ws.Range("N19").Value = 17.256637168141

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.: WTD 2023 count (C) flips from the "0" placeholder to a
# real count of 1
# ---------------------------------------------------------------------------
Set-NumberCell $ws "C20" "D20" 1
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = -76.923076923076
$ws.Range("I20").Value = 73
$ws.Range("J20").Value = 116
$ws.Range("K20").Value = -37.068965517241
$ws.Range("L20").Value = 32.727272727272
$ws.Range("M20").Value = 102.777777777778
$ws.Range("N20").Value = -89.006024096385

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = -25
$ws.Range("F21").Value = 42
$ws.Range("G21").Value = 45
$ws.Range("H21").Value = -6.666666666666
$ws.Range("I21").Value = 496
$ws.Range("J21").Value = 476
$ws.Range("K21").Value = 4.201680672268
$ws.Range("L21").Value = 66.442953020134
$ws.Range("M21").Value = 44.186046511627
$ws.Range("N21").Value = -63.927272727272

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 41.666666666666
$ws.Range("F24").Value = 47
$ws.Range("G24").Value = 56
$ws.Range("H24").Value = -16.071428571428
$ws.Range("I24").Value = 462
$ws.Range("J24").Value = 457
$ws.Range("K24").Value = 1.094091903719
$ws.Range("L24").Value = 81.889763779527
$ws.Range("M24").Value = -8.151093439363

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 14
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = -17.647058823529
$ws.Range("I25").Value = 180
$ws.Range("J25").Value = 173
$ws.Range("K25").Value = 4.046242774566
$ws.Range("L25").Value = 35.338345864661
$ws.Range("M25").Value = -10.447761194029

# ---------------------------------------------------------------------------
# Row 26 - UCR Rape*: 28-day 2022 count (G) & %chg (H) go to text placeholders
# ---------------------------------------------------------------------------
Set-TextCell $ws "G26" "F26" "0"
Set-TextCell $ws "H26" "F26" "***.*"

# ---------------------------------------------------------------------------
# Row 27 - Other Sex Crimes: WTD 2023 (C) and 28-day 2023 (F) flip from "0"
# placeholders to real counts of 1
# ---------------------------------------------------------------------------
Set-NumberCell $ws "C27" "G27" 1
Set-NumberCell $ws "F27" "G27" 1
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 17
$ws.Range("K27").Value = -5.555555555555
$ws.Range("L27").Value = 41.666666666666

# ---------------------------------------------------------------------------
# Row 30 - Hate Crimes: 28-day 2023 (F) and YTD 2023 (I) flip from "0" to 1,
# while 28-day 2022 (G) and %chg (H) flip from real numbers to placeholders
# ---------------------------------------------------------------------------
Set-NumberCell $ws "F30" "J30" 1
Set-TextCell $ws "G30" "E30" "0"
Set-TextCell $ws "H30" "E30" "***.*"
Set-NumberCell $ws "I30" "J30" 1
$ws.Range("K30").Value = -83.333333333333
$ws.Range("L30").Value = -50

Write-Output "edits applied"
